# Add three new Venue/County rows to the bottom of the venue list
# (matches the rows added in the source OOXML diff: rows 125-127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A125").Value = "Find Insurance Celtic Park, Derry"
$ws.Range("B125").Value = "Derry"

$ws.Range("A126").Value = "VBC Cloghan, Castleblayney"
$ws.Range("B126").Value = "Monaghan"

$ws.Range("A127").Value = "Cappoquin Logistics Fraher Field"
$ws.Range("B127").Value = "Waterford"

# Restore the selection the author left the sheet on after the edit.
$ws.Range("A114").Select()
